$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.040.77'
$ws.Range('E2').Value = '  +6.00%  '

$ws.Range('D3').Value = '3.118.36'
$ws.Range('E3').Value = '  +3.76%  '

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').Value = '''587.92'
$ws.Range('E5').Value = '  +4.29%  '

$ws.Range('D6').Value = '''143.91'
$ws.Range('E6').Value = '  +3.27%  '

$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').Value = '3.104.71'
$ws.Range('E8').Value = '  +3.73%  '

$ws.Range('E9').Value = '  +2.38%  '

$ws.Range('D10').Value = '''0.146'
$ws.Range('E10').Value = '  +9.68%  '

$ws.Range('D11').Value = '''5.75'
$ws.Range('E11').Value = '  +10.16%  '

$ws.Range('D12').Value = '''0.470'
$ws.Range('E12').Value = '  +2.27%  '

$ws.Range('E13').Value = '  +5.32%  '

$ws.Range('E14').Value = '  +5.13%  '

$ws.Range('E15').Value = '  +0.66%  '

$ws.Range('D16').Value = '3.631.86'
$ws.Range('E16').Value = '  +3.82%  '

$ws.Range('E17').Value = '  -1.21%  '

$ws.Range('D18').Value = '62.995.83'
$ws.Range('E18').Value = '  +6.14%  '

$ws.Range('D19').Value = '3.113.28'
$ws.Range('E19').Value = '  +3.83%  '

$ws.Range('D20').Value = '''453.91'
$ws.Range('E20').Value = '  +5.02%  '

$ws.Range('D21').Value = '''14.10'
$ws.Range('E21').Value = '  +3.22%  '

$ws.Range('D22').Value = '''0.735'
$ws.Range('E22').Value = '  +1.50%  '

$ws.Range('D23').Value = '''7.55'
$ws.Range('E23').Value = '  +5.49%  '

$ws.Range('D24').Value = '''13.63'
$ws.Range('E24').Value = '  +0.78%  '

$ws.Range('D25').Value = '''82.08'
$ws.Range('E25').Value = '  +1.85%  '

$ws.Range('E26').Value = '  +0.19%  '

$ws.Range('D27').Value = '''2.27'
$ws.Range('E27').Value = '  +1.08%  '

$ws.Range('E28').Value = '  +5.79%  '

$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '''8.29'
$ws.Range('E29').Value = '  +4.59%  '

$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  +0.13%  '

$ws.Range('D31').Value = '''6.86'
$ws.Range('E31').Value = '  +11.56%  '

$ws.Range('E32').Value = '  +11.75%  '

$ws.Range('D33').Value = '''27.15'
$ws.Range('E33').Value = '  +5.05%  '

$ws.Range('D34').Value = '''1.05'
$ws.Range('E34').Value = '  +4.12%  '

$ws.Range('D35').Value = '0.0₃0808'
$ws.Range('E35').Value = '  +6.01%  '

$ws.Range('E36').Value = '  +8.60%  '

$ws.Range('E37').Value = '  +1.21%  '

$ws.Range('E38').Value = '  +3.72%  '

$ws.Range('D39').Value = '''3.04'
$ws.Range('E39').Value = '  +9.83%  '

$ws.Range('D40').Value = '''8.81'
$ws.Range('E40').Value = '  +1.27%  '

$ws.Range('D41').Value = '''428.52'
$ws.Range('E41').Value = '  +3.73%  '

$ws.Range('D42').Value = '2.958.99'
$ws.Range('E42').Value = '  +6.46%  '

$ws.Range('E43').Value = '  +5.71%  '

$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = '''0.112'
$ws.Range('E44').Value = '  +3.20%  '

$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '''0.276'
$ws.Range('E45').Value = '  +8.89%  '

$ws.Range('E46').Value = '  +7.35%  '

$ws.Range('D47').Value = '''125.66'
$ws.Range('E47').Value = '  +1.57%  '

$ws.Range('D49').Value = '''34.74'
$ws.Range('E49').Value = '  -0.65%  '

$ws.Range('E50').Value = '  +1.00%  '

$ws.Range('D51').Value = '''24.82'
$ws.Range('E51').Value = '  +5.10%  '
